$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand Table1 to cover the new rows (A1:B296)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B296"))

# Populate the newly added lookup rows
$ws.Range("A259").Value = "(9S1..)"
$ws.Range("B259").Value = "White"
$ws.Range("A260").Value = "(9S3..)"
$ws.Range("B260").Value = "Black or Black British"
$ws.Range("A261").Value = "(XaFwD)"
$ws.Range("B261").Value = "White"
$ws.Range("A262").Value = "(XaFwE)"
$ws.Range("B262").Value = "White"
$ws.Range("A263").Value = "(XaFwH)"
$ws.Range("B263").Value = "Black or Black British"
$ws.Range("A264").Value = "(XaFwz)"
$ws.Range("B264").Value = "Asian or Asian British"
$ws.Range("A265").Value = "(XaFx1)"
$ws.Range("B265").Value = "Other Ethnic Group"
$ws.Range("A266").Value = "(XaIuh)"
$ws.Range("B266").Value = "White"
$ws.Range("A267").Value = "(XaJQv)"
$ws.Range("B267").Value = "Mixed"
$ws.Range("A268").Value = "(XaJQv) British or mixed British - ethnic category"
$ws.Range("B268").Value = "Mixed"
$ws.Range("A269").Value = "(XaJQw)"
$ws.Range("B269").Value = "Mixed"
$ws.Range("A270").Value = "(XaJQx)"
$ws.Range("B270").Value = "White"
$ws.Range("A271").Value = "(XaJQy)"
$ws.Range("B271").Value = "White"
$ws.Range("A272").Value = "(XaJR6)"
$ws.Range("B272").Value = "Black or Black British"
$ws.Range("A273").Value = "(XaJR7)"
$ws.Range("B273").Value = "Black or Black British"
$ws.Range("A274").Value = "(XaJRA)"
$ws.Range("B274").Value = "Black or Black British"
$ws.Range("A275").Value = "(XaJRB)"
$ws.Range("B275").Value = "Not stated"
$ws.Range("A276").Value = "(XaJRC)"
$ws.Range("B276").Value = "White"
$ws.Range("A277").Value = "(XaJRD)"
$ws.Range("B277").Value = "White"
$ws.Range("A278").Value = "(XaJSD)"
$ws.Range("B278").Value = "Other Ethnic Group"
$ws.Range("A279").Value = "(XaJSE)"
$ws.Range("B279").Value = "Other Ethnic Group"
$ws.Range("A280").Value = "(XaJSF)"
$ws.Range("B280").Value = "Other Ethnic Group"
$ws.Range("A281").Value = "(XaJSG)"
$ws.Range("B281").Value = "Other Ethnic Group"
$ws.Range("A282").Value = "(XaJSK)"
$ws.Range("B282").Value = "Other Ethnic Group"
$ws.Range("A283").Value = "(XaJSa)"
$ws.Range("B283").Value = "Other Ethnic Group"
$ws.Range("A284").Value = "(XaQEa)"
$ws.Range("B284").Value = "White"
$ws.Range("A285").Value = "(XaR4o)"
$ws.Range("B285").Value = "Other Ethnic Group"
$ws.Range("A286").Value = "(XactH)"
$ws.Range("B286").Value = "White"
$ws.Range("A287").Value = "(XactI)"
$ws.Range("B287").Value = "Black or Black British"
$ws.Range("A288").Value = "(XactK)"
$ws.Range("B288").Value = "White"
$ws.Range("A289").Value = "(Xactk)"
$ws.Range("B289").Value = "Asian or Asian British"
$ws.Range("A290").Value = "(Xacut)"
$ws.Range("B290").Value = "White"
$ws.Range("A291").Value = "British or Mixed British"
$ws.Range("B291").Value = "Mixed"
$ws.Range("A292").Value = "Other Asian"
$ws.Range("B292").Value = "Asian or Asian British"
$ws.Range("A293").Value = "Other Mixed"
$ws.Range("B293").Value = "Mixed"
$ws.Range("A294").Value = "White & Asian"
$ws.Range("B294").Value = "Mixed"
$ws.Range("A295").Value = "White & Black African"
$ws.Range("B295").Value = "Mixed"
$ws.Range("A296").Value = "White & Black Caribbean"
$ws.Range("B296").Value = "Mixed"

# Reset the active selection to A2 (matches post-edit saved state)
$ws.Range("A2").Select()
